$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: after a paragraph's existing text (found via $searchText), append
# a new run containing a single space (plain formatting) followed by a new
# italic run containing $word. A temporary bookmark is dropped at the
# original/space boundary and immediately removed again purely to force the
# engine to keep that boundary as two distinct <w:r> elements instead of
# silently re-coalescing them (mirrors the way Word itself leaves behind
# separate runs when content is composed across more than one edit).
# ---------------------------------------------------------------------------
function Add-ItalicSuffix($doc, $searchText, $word) {
    $rng = $doc.Content
    $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $endPos = $rng.End

    $spaceRng = $doc.Range($endPos, $endPos)
    $spaceRng.InsertAfter(" ")

    $wordRng = $doc.Range($spaceRng.End, $spaceRng.End)
    $wordRng.InsertAfter($word)
    $wordRng.Font.Italic = 1

    $splitRng = $doc.Range($endPos, $endPos)
    $doc.Bookmarks.Add("TempSplit", $splitRng)
    $doc.Bookmarks("TempSplit").Delete()
}

# ---------------------------------------------------------------------------
# Helper: force a run boundary at a given document position without changing
# any text, by dropping then immediately deleting a scratch bookmark there.
# ---------------------------------------------------------------------------
function Split-RunAt($doc, $pos) {
    $r = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add("TempSplit", $r)
    $doc.Bookmarks("TempSplit").Delete()
}

# 1) Investments indices -> + " fake" (italic)
Add-ItalicSuffix $d "Investments indices" "fake"

# 2) Stock prices (S&P/TSX/NYSE/NASDAQ) -> + " real" (italic)
Add-ItalicSuffix $d "Stock prices (S&P/TSX/NYSE/NASDAQ)" "real"

# 3) S&P/TSX Composite Index -> + " real" (italic)
Add-ItalicSuffix $d "S&P/TSX Composite Index" "real"

# 4) ...ajor indices and commodities -> + " fake" (italic)
Add-ItalicSuffix $d "ajor indices and commodities" "fake"

# 5) Quotes -> + " fake" (italic)
Add-ItalicSuffix $d "Quotes" "fake"

# 6) Mutual Funds (active/passive) -> + " fake" (italic)
Add-ItalicSuffix $d " (active/passive)" "fake"

# ---------------------------------------------------------------------------
# 7) "Custom GUI design" -> "Custom " + _GoBack bookmark + "graphics design"
#    split into the run layout: "Custom" | " " | [_GoBack] | "graphics" | " design"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Custom GUI design", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pStart = $rng.Start

# Split "Custom" | " GUI design"
Split-RunAt $d ($pStart + 6)

# Move _GoBack to sit between " " and "GUI" (this also removes it from its
# old location at the end of the document, since _GoBack is a singleton).
$goBackPos = $pStart + 7
$goBackRng = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRng)

# Replace "GUI" with "graphics"
$guiRng = $d.Content
$guiRng.Find.Execute("GUI", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$guiRng.Text = "graphics"

# Split "graphics" | " design"
Split-RunAt $d $guiRng.End

Write-Host "All edits applied"
